$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.011.80'
$ws.Range("E2").Value = '  +4.70%  '
$ws.Range("D3").Value = '3.415.61'
$ws.Range("E3").Value = '  +3.44%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.27'
$ws.Range("E5").Value = '  +3.66%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.88'
$ws.Range("E6").Value = '  +9.71%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '3.414.70'
$ws.Range("E8").Value = '  +3.38%  '
$ws.Range("E9").Value = '  +2.45%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.50'
$ws.Range("E10").Value = '  +2.31%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.128'
$ws.Range("E11").Value = '  +10.45%  '
$ws.Range("E12").Value = '  +7.16%  '
$ws.Range("D13").Value = '4.001.20'
$ws.Range("E13").Value = '  +3.32%  '
$ws.Range("E14").Value = '  +1.96%  '
$ws.Range("E15").Value = '  +9.11%  '
$ws.Range("D16").Value = '3.410.08'
$ws.Range("E16").Value = '  +3.09%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '25.50'
$ws.Range("E17").Value = '  +7.24%  '
$ws.Range("D18").Value = '62.008.68'
$ws.Range("E18").Value = '  +4.25%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.20'
$ws.Range("E19").Value = '  +7.97%  '
$ws.Range("E20").Value = '  +5.68%  '
$ws.Range("E21").Value = '  +8.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '390.25'
$ws.Range("E22").Value = '  +12.35%  '
$ws.Range("E23").Value = '  +4.27%  '
$ws.Range("D24").Value = '3.552.65'
$ws.Range("E24").Value = '  +3.40%  '
$ws.Range("E25").Value = '  +19.28%  '
$ws.Range("E26").Value = '  +0.13%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '71.58'
$ws.Range("E27").Value = '  +4.78%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.60'
$ws.Range("E28").Value = '  +12.01%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.68'
$ws.Range("E29").Value = '  +6.41%  '
$ws.Range("E30").Value = '  -0.16%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.31'
$ws.Range("E31").Value = '  +7.58%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.160'
$ws.Range("E32").Value = '  +7.00%  '
$ws.Range("E33").Value = '  +5.18%  '
$ws.Range("D34").Value = '3.447.33'
$ws.Range("E34").Value = '  +3.47%  '
$ws.Range("E36").Value = '  +4.47%  '
$ws.Range("E37").Value = '  +5.95%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.02'
$ws.Range("E38").Value = '  +4.40%  '
$ws.Range("E39").Value = '  +7.47%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '162.87'
$ws.Range("E40").Value = '  +3.17%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0793'
$ws.Range("E41").Value = '  +7.16%  '
$ws.Range("E42").Value = '  +16.66%  '
$ws.Range("E43").Value = '  +7.15%  '
$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.999'
$ws.Range("E44").Value = '  -0.07%  '
$ws.Range("B45").Value = 'ONDO'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.23'
$ws.Range("E45").Value = '  +6.78%  '
$ws.Range("E46").Value = '  +5.27%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '25.11'
$ws.Range("E47").Value = '  +11.51%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '41.76'
$ws.Range("E48").Value = '  +3.63%  '
$ws.Range("E49").Value = '  +4.73%  '
$ws.Range("E50").Value = '  +8.09%  '
$ws.Range("D51").Value = '2.377.91'
$ws.Range("E51").Value = '  +10.61%  '
